$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 155, shifting existing rows 155..216 down to 156..217
$ws.Rows.Item(155).Insert()

# Populate the newly inserted row 155 with the new data record
$ws.Cells.Item(155, 1).Value  = 9
$ws.Cells.Item(155, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(155, 3).Value  = "Metropolitana"
$ws.Cells.Item(155, 4).Value  = 44468
$ws.Cells.Item(155, 5).Value  = 13
$ws.Cells.Item(155, 6).Value  = 100112044
$ws.Cells.Item(155, 7).Value  = "Perejil"
$ws.Cells.Item(155, 8).Value  = "Sin especificar"
$ws.Cells.Item(155, 9).Value  = "Primera"
$ws.Cells.Item(155, 10).Value = 106
$ws.Cells.Item(155, 11).Value = 10000
$ws.Cells.Item(155, 12).Value = 12000
$ws.Cells.Item(155, 13).Value = 11000
$ws.Cells.Item(155, 14).Value = "`$/docena de atados"
$ws.Cells.Item(155, 15).Value = "Región Metropolitana"
$ws.Cells.Item(155, 16).Value = 3667
$ws.Cells.Item(155, 17).Value = 3
$ws.Cells.Item(155, 18).Value = "Hortaliza"
